$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 162, shifting the existing rows 162:239 down to 163:240.
$ws.Rows("162:162").Insert()

# Populate the newly inserted row 162 with the new data record.
$ws.Range("A162").Value = 6
$ws.Range("B162").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C162").Value = "Metropolitana"
$ws.Range("D162").Value = 44784
$ws.Range("E162").Value = 13
$ws.Range("F162").Value = 100112022
$ws.Range("G162").Value = "Arveja Verde"
$ws.Range("H162").Value = "Perfection"
$ws.Range("I162").Value = "Primera"
$ws.Range("J162").Value = 220
$ws.Range("K162").Value = 38000
$ws.Range("L162").Value = 40000
$ws.Range("M162").Value = 39091
$ws.Range("N162").Value = "`$/saco 25 kilos"
$ws.Range("O162").Value = "Provincia de Huasco"
$ws.Range("P162").Value = 1564
$ws.Range("Q162").Value = 25
$ws.Range("R162").Value = "Hortaliza"
